# MathuraCorreo/correos.xlsx - "Added some examples in email sendind"
#
# Replace the email address in A3 (and its mailto hyperlink) and add a new
# row 4 with another example email address + hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The COM hyperlink shim only supports a sheet-wide Hyperlinks.Delete(), so
# drop every hyperlink on the sheet and rebuild all three (A2 unchanged,
# A3 updated, A4 new) in order so the relationship ids line up again.
$ws.Range("A2").Hyperlinks.Delete()

$ws.Range("A3").Value = "luisfervillaalta@gmail.com"
$ws.Range("A4").Value = "jesssortigoza@gmail.com"

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:elvisavfc65@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:luisfervillaalta@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:jesssortigoza@gmail.com")

# Hyperlinks.Add applies its own ad-hoc style; reapply the workbook's
# built-in hyperlink cell style so A2:A4 keep matching formatting.
$ws.Range("A2:A4").Style = "Hipervínculo"

$ws.Range("A5").Select() | Out-Null
